$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the symmetric waypoint references (Xn) into the "connected-to" list
# of the nodes that those waypoints connect, so connections become symmetric.
$ws.Range("F7").Value()  = "7;X1"
$ws.Range("F15").Value() = "15;X1"
$ws.Range("F9").Value()  = "7;11;X2"
$ws.Range("F17").Value() = "15;17;X2"
$ws.Range("F13").Value() = "10;11;X3"
$ws.Range("F19").Value() = "17;22;X3"
$ws.Range("F24").Value() = "19;20;24;X4"
$ws.Range("F30").Value() = "28;30;31;X4"
$ws.Range("F27").Value() = "22;25;X5"
$ws.Range("F32").Value() = "29;32;X5"

# Update the active selection to match the saved workbook state
$ws.Range("F32").Select()
